$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Insert a new column at I (shifts old I..O -> J..P) for the new
# "Seller Signatory Emails" field, mirroring the existing "Email *"
# column (D) in content/format/hyperlinks.
# -----------------------------------------------------------------
$ws.Columns("I:I").Insert()

# Match the new column's width to its left neighbour (Pan, col H)
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth

# Copy values + formatting from the Email column (D) into the new column (I)
$ws.Range("D1:D7").Copy($ws.Range("I1:I7"))

# Set the new header text
$ws.Range("I1").Value2 = "Seller Signatory Emails"

# Re-create hyperlinks on the copied e-mail cells (mirrors D4/D6/D7)
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:emp3@myfirm.com")
$ws.Hyperlinks.Add($ws.Range("I6"), "mailto:emp1@investor1.com")
$ws.Hyperlinks.Add($ws.Range("I7"), "mailto:emp1@investor2.com")

# Hyperlinks.Add() re-applies its own format xf; restore the plain
# "Hyperlink" cell style so it matches column D's existing xf entry
# instead of creating a near-duplicate style.
$ws.Range("I4").Style = "Hyperlink"
$ws.Range("I6").Style = "Hyperlink"
$ws.Range("I7").Style = "Hyperlink"

# -----------------------------------------------------------------
# The previously-unused "applyFill" style variant is no longer
# referenced; normalise the "Update Only / DP / Client Id" columns
# (now N:P, previously M:O) back onto the plain "Normal 2" style.
# -----------------------------------------------------------------
$ws.Range("N1:P1").Style = "Normal 2"
$ws.Range("N2:N7").Style = "Normal 2"

# -----------------------------------------------------------------
# Update the sheet view: scroll so column D is left-most, and select
# the freshly added column's data range.
# -----------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("I2:I7").Select()
